# Refresh the "Price" column (D) of the crypto symbol list with the latest
# values from the scrape run on 2022-12-15 10:30 UTC.
#
# The sheet stores every "numeric" column as text (so formatting like
# "0.03405" with trailing/leading zeros round-trips exactly), so a plain
# `.Value = "..."` assignment would let Excel's type-inference turn the
# numeric-looking string into a real number. To avoid that we briefly mark
# the cell as Text ("@") before writing the new value, then restore the
# default "Normal" style so we don't leave any stray number-formatting
# behind on cells that originally had none.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (in column D) -> new Price text, taken from the refreshed feed.
$priceUpdates = [ordered]@{
    2  = "265.24"
    3  = "22.88"
    4  = "6.206"
    5  = "0.06162"
    6  = "3.561"
    7  = "6.700"
    8  = "1.361"
    9  = "0.8126"
    10 = "0.1589"
    11 = "0.08192"
    12 = "0.03394"
    13 = "0.03151"
    14 = "0.09229"
    15 = "3.897"
    16 = "0.001714"
    17 = "0.04840"
    18 = "0.0006261"
    19 = "0.006177"
    20 = "0.006266"
    21 = "0.001097"
    22 = "0.0001500"
    23 = "3.699"
    25 = "0.3383"
    27 = "0.0002682"
    40 = "0.04589"
    41 = "0.006990"
    42 = "0.1132"
    43 = "0.003400"
    44 = "0.01065"
    45 = "0.00006104"
    47 = "0.7701"
    48 = "0.1885"
    49 = "0.00002100"
    50 = "0.01240"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}
